$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.050.61"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.321.01"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.40"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.20"
$ws.Range("E6").Value = "  +5.52%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.320.68"
$ws.Range("E8").Value = "  +6.32%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("E11").Value = "  +6.20%  "
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "3.867.11"
$ws.Range("E15").Value = "  +6.10%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "3.323.62"
$ws.Range("E17").Value = "  +6.21%  "
$ws.Range("D18").Value = "64.118.24"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.91"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.06"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +5.88%  "
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.72"
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.12"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").Value = "  +4.41%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.72"
$ws.Range("E31").Value = "  +10.95%  "
$ws.Range("E32").Value = "  +5.91%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "0.0₃0767"
$ws.Range("E37").Value = "  +7.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.34"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0404"
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "436.48"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").Value = "3.061.11"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.81"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.46"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("E47").Value = "  +3.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.24"
$ws.Range("E48").Value = "  +14.14%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  +1.75%  "
